$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3668.6428
$ws.Range("I43").Value = 2755.1667
$ws.Range("J43").Value = 4353.75
$ws.Range("K43").Value = 2755.1667
$ws.Range("L43").Value = 4353.75
$ws.Range("M43").Value = -2686.1667
$ws.Range("N43").Value = -4491.75
$ws.Range("H46").Value = 2800
$ws.Range("J46").Value = 2800
$ws.Range("L46").Value = 8400
$ws.Range("N46").Value = -8638
$ws.Range("H51").Value = 10420770
$ws.Range("J51").Value = 4279.5
$ws.Range("L51").Value = 4279.5
$ws.Range("N51").Value = -5247.5
$ws.Range("H60").Value = 2800
$ws.Range("J60").Value = 2800
$ws.Range("L60").Value = 8400
$ws.Range("N60").Value = -9368
$ws.Range("H74").Value = 2201.5
$ws.Range("I74").Value = 2201.5
$ws.Range("K74").Value = 2201.5
$ws.Range("M74").Value = -1265.5
$ws.Range("H77").Value = 2201.5
$ws.Range("I77").Value = 2201.5
$ws.Range("K77").Value = 11007.5
$ws.Range("M77").Value = -6327.5
$ws.Range("H98").Value = 2133.8572
$ws.Range("I98").Value = 1144.3077
$ws.Range("K98").Value = 1144.3077
$ws.Range("M98").Value = 353.6922999999999
$ws.Range("H111").Value = 1131.4667
$ws.Range("I111").Value = 985.75
$ws.Range("K111").Value = 2957.25
$ws.Range("M111").Value = 109.75
$ws.Range("H122").Value = 2133.8572
$ws.Range("I122").Value = 1144.3077
$ws.Range("K122").Value = 3432.9231
$ws.Range("M122").Value = -982.9231
$ws.Range("H137").Value = 14201.68
$ws.Range("I137").Value = 24198.54
$ws.Range("K137").Value = 72595.62
$ws.Range("M137").Value = -70045.62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 601000.8
$ws.Range("I74").Value = 601000.8
$ws.Range("K74").Value = 601000.8
$ws.Range("M74").Value = -600126.8
$ws.Range("H77").Value = 601000.8
$ws.Range("I77").Value = 601000.8
$ws.Range("K77").Value = 3005004
$ws.Range("M77").Value = -3000636
$ws.Range("H97").Value = 1321.1459
$ws.Range("I97").Value = 923.55884
$ws.Range("J97").Value = 2286.7144
$ws.Range("K97").Value = 923.55884
$ws.Range("L97").Value = 2286.7144
$ws.Range("M97").Value = -427.55884
$ws.Range("N97").Value = -3278.7144
$ws.Range("H122").Value = 2543.1667
$ws.Range("I122").Value = 2528.3928
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 7585.178400000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -5135.178400000001
$ws.Range("N122").Value = -13150
$ws.Range("H132").Value = 1342.6305
$ws.Range("I132").Value = 1012.9643
$ws.Range("J132").Value = 1855.4445
$ws.Range("K132").Value = 3038.8929
$ws.Range("L132").Value = 5566.333500000001
$ws.Range("M132").Value = -508.8928999999998
$ws.Range("N132").Value = -10626.3335
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3735.6
$ws.Range("I105").Value = 2799.25
$ws.Range("K105").Value = 2799.25
$ws.Range("M105").Value = -1052.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8212.620999999999
$ws.Range("J99").Value = 15347.4
$ws.Range("L99").Value = 15347.4
$ws.Range("N99").Value = -18343.4
$ws.Range("H107").Value = 484.36667
$ws.Range("I107").Value = 318.73914
$ws.Range("J107").Value = 1028.5714
$ws.Range("K107").Value = 318.73914
$ws.Range("L107").Value = 1028.5714
$ws.Range("M107").Value = 1601.26086
$ws.Range("N107").Value = -4868.5714
$ws.Range("H126").Value = 8212.620999999999
$ws.Range("J126").Value = 15347.4
$ws.Range("L126").Value = 46042.2
$ws.Range("N126").Value = -50982.2
$ws.Range("H132").Value = 18383.428
$ws.Range("I132").Value = 19891.562
$ws.Range("K132").Value = 59674.686
$ws.Range("M132").Value = -57144.686
$ws.Range("H138").Value = 84188.64
$ws.Range("J138").Value = 84188.64
$ws.Range("L138").Value = 84188.64
$ws.Range("N138").Value = -94468.64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 817178.3
$ws.Range("I2").Value = 926023.8
$ws.Range("K2").Value = 5556142.800000001
$ws.Range("M2").Value = -5556029.800000001
$ws.Range("H12").Value = 133
$ws.Range("J12").Value = 114.44444
$ws.Range("L12").Value = 343.33332
$ws.Range("N12").Value = -689.33332
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H55").Value = 2542.7144
$ws.Range("J55").Value = 3499.75
$ws.Range("L55").Value = 10499.25
$ws.Range("N55").Value = -10853.25
$ws.Range("H68").Value = 4593.1724
$ws.Range("J68").Value = 4884.5386
$ws.Range("L68").Value = 14653.6158
$ws.Range("N68").Value = -16275.6158
$ws.Range("H71").Value = 4593.1724
$ws.Range("J71").Value = 4884.5386
$ws.Range("L71").Value = 43960.8474
$ws.Range("N71").Value = -52072.8474
$ws.Range("H80").Value = 4884.3335
$ws.Range("I80").Value = 4500
$ws.Range("K80").Value = 13500
$ws.Range("M80").Value = -12564
$ws.Range("H83").Value = 4884.3335
$ws.Range("I83").Value = 4500
$ws.Range("K83").Value = 40500
$ws.Range("M83").Value = -35820
$ws.Range("H116").Value = 9195.608
$ws.Range("I116").Value = 3833
$ws.Range("K116").Value = 11499
$ws.Range("M116").Value = -8057
$ws.Range("H132").Value = 1635.591
$ws.Range("I132").Value = 1565.8182
$ws.Range("J132").Value = 1705.3636
$ws.Range("K132").Value = 14092.3638
$ws.Range("L132").Value = 15348.2724
$ws.Range("M132").Value = -11562.3638
$ws.Range("N132").Value = -20408.2724
$ws.Range("H136").Value = 5910.5
$ws.Range("I136").Value = 2287.6667
$ws.Range("K136").Value = 6863.000100000001
$ws.Range("M136").Value = -1763.000100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 704
$ws.Range("I2").Value = 356
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 356
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -243
$ws.Range("N2").Value = -1626
$ws.Range("H102").Value = 2942.8
$ws.Range("I102").Value = 2428.8125
$ws.Range("K102").Value = 2428.8125
$ws.Range("M102").Value = -806.8125
$ws.Range("H107").Value = 618.3333
$ws.Range("J107").Value = 618.3333
$ws.Range("L107").Value = 618.3333
$ws.Range("N107").Value = -4458.3333
$ws.Range("H126").Value = 2468.3157
$ws.Range("I126").Value = 2182.3333
$ws.Range("J126").Value = 2958.5715
$ws.Range("K126").Value = 6546.999899999999
$ws.Range("L126").Value = 8875.7145
$ws.Range("M126").Value = -4076.999899999999
$ws.Range("N126").Value = -13815.7145
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1586.8823
$ws.Range("I61").Value = 1644.8
$ws.Range("J61").Value = 1152.5
$ws.Range("K61").Value = 1644.8
$ws.Range("L61").Value = 1152.5
$ws.Range("M61").Value = -1442.8
$ws.Range("N61").Value = -1556.5
$ws.Range("H113").Value = 1586.8823
$ws.Range("I113").Value = 1644.8
$ws.Range("J113").Value = 1152.5
$ws.Range("K113").Value = 1644.8
$ws.Range("L113").Value = 1152.5
$ws.Range("M113").Value = 525.2
$ws.Range("N113").Value = -5492.5
$ws.Range("H132").Value = 2200.6667
$ws.Range("I132").Value = 2301
$ws.Range("K132").Value = 6903
$ws.Range("M132").Value = -4373
$ws.Range("H134").Value = 124962.336
$ws.Range("I134").Value = 124962
$ws.Range("J134").Value = 124963
$ws.Range("K134").Value = 124962
$ws.Range("L134").Value = 124963
$ws.Range("N134").Value = -135103
$ws.Range("M134").Value = -119892
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 738.5238000000001
$ws.Range("I113").Value = 586.3333
$ws.Range("K113").Value = 1758.9999
$ws.Range("M113").Value = 411.0001
$ws.Range("H132").Value = 3466.5625
$ws.Range("I132").Value = 3497.5715
$ws.Range("K132").Value = 10492.7145
$ws.Range("M132").Value = -7962.7145
$ws.Range("H136").Value = 12714.426
$ws.Range("I136").Value = 13344.804
$ws.Range("K136").Value = 40034.412
$ws.Range("M136").Value = -37484.412
